$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update row 11 (second table: Date/Time/Hours/Activity in O:S)
$ws.Range("S11").Value = "Redoing all previous tests. Didn" + [char]0x2019 + "t realize vsync was on... Performed graphics quality test. Performed lighting test."
$ws.Range("P11").Value = "8.15 - 15.15"
$ws.Range("R11").Value = 7

# Add new row 27 of data (first table: Date/Time/Hours/Activity in A:E)
$ws.Range("A27").Value = (Get-Date -Year 2022 -Month 7 -Day 5).Date
$ws.Range("B27").Value = "17.00 - 20.30"
$ws.Range("D27").Value = 3.5
$ws.Range("E27").Value = "Finalized game. Preparing handin."

# Update selection to match the recorded cursor position
$ws.Range("J33").Select()
